# Socket client application doc - apply edits described by the unified diff.
# Strategy: locate each target paragraph (or paragraph range) by its current
# text, then replace the whole paragraph(s) with freshly authored OOXML via
# Range.InsertXML. This merges runs, drops stale <w:proofErr/> marks, and
# (for the big restructuring) reorders/creates paragraphs and relocates the
# _GoBack bookmark, all in one atomic operation per hunk.

$d = $word.ActiveDocument

function Find-ParaByText($doc, [string]$matchText, [bool]$exact) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($exact) {
            if ($t -eq $matchText) { return $p }
        } else {
            if ($t.StartsWith($matchText)) { return $p }
        }
    }
    return $null
}

# --- Hunk 1: merge "Following are the technologies and apis being used for this project," run, drop proofErr ---
$para1 = Find-ParaByText $d "Following are the technologies and" $false
$range1 = $d.Range($para1.Range.Start, $para1.Range.End - 1)
$range1.InsertXML('<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Following are the technologies and apis being used for this project,</w:t></w:r></w:p>')

# --- Hunk 2: merge "Jaxb-api , jaxb-impl jars " run, drop proofErr ---
$para2 = Find-ParaByText $d "Jaxb-api" $false
$range2 = $d.Range($para2.Range.Start, $para2.Range.End - 1)
$range2.InsertXML('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Jaxb-api , jaxb-impl jars </w:t></w:r></w:p>')

# --- Hunk 3: reorder paragraphs after "JUnit", insert two new paragraphs, merge
#     "Junit is used..." run, drop proofErr marks, and move the _GoBack
#     bookmark to the new trailing paragraph ---
$startPara = Find-ParaByText $d "JUnit" $true
$endPara = Find-ParaByText $d "Junit is used to provide the unit test for the different business case scenarios. " $true
$range3 = $d.Range($startPara.Range.Start, $endPara.Range.End)
$range3.InsertXML('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>JUnit</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>This application provides the socket client functionality using which we can write an xml message to a socket output stream and read the response from socket’s input stream</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>This solution is implemented using observer pattern where a Socket client observer observes the successful writing of xml message to the socket’s input stream, once the xml message sent the observer is notified accordingly.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The implementation is abstracted behind the </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>SocketSendAndReceive</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> interface, the reason is to provide extensibility in future.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Junit is used to provide the unit test for the different business case scenarios. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The application is testable via a main class or using the unit test case in </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>SocketClientTest</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> class</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p>')

# --- Hunk 4: merge the three maven-command paragraphs into single runs, drop proofErr ---
$para4a = Find-ParaByText $d "To run unit " $false
$para4c = Find-ParaByText $d "To build and install" $false
$range4 = $d.Range($para4a.Range.Start, $para4c.Range.End - 1)
$range4.InsertXML('<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>To run unit test : mvn clean test</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>To Build and package : mvn clean package</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>To build and install : mvn clean install.</w:t></w:r></w:p>')

# --- Hunk 5: merge "README.txt , which is under root of the project." run, drop proofErr ---
$para5 = Find-ParaByText $d "The final result will be written to a text file called" $false
$range5 = $d.Range($para5.Range.Start, $para5.Range.End - 1)
$range5.InsertXML('<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The final result will be written to a text file called </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>README.txt , which is under root of the project.</w:t></w:r></w:p>')

Write-Output "All edits applied"
